$p = $ppt.ActivePresentation
$s = $p.Slides.Item(6)
$shape = $s.Shapes.Item(2)
$tbl = $shape.Table
$tbl.ApplyStyle("{5DA735B5-EC3B-4E9F-AEA3-90FDA6914C37}")
